$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above the current row 94, shifting existing rows 94:147 down to 98:151.
$ws.Range("A94:A97").EntireRow.Insert()

# Common (boilerplate) values shared by every data row in this sheet.
$mercadoId = 4
$mercado   = "Feria Lagunitas de Puerto Montt"
$region    = "Los Lagos"
$codreg    = 10
$tipo      = "Fruta"
$productoId = 100103
$producto   = "Frutos de hueso (carozo)"
$categoriaId = 100103002
$categoria   = "Ciruela"
$origenOHiggins = "Región de O'Higgins"

# New row 94: Black Amber, Primera
$ws.Cells.Item(94, 1).Value = $mercadoId
$ws.Cells.Item(94, 2).Value = $mercado
$ws.Cells.Item(94, 3).Value = $region
$ws.Cells.Item(94, 4).Value = 44603
$ws.Cells.Item(94, 5).Value = $codreg
$ws.Cells.Item(94, 6).Value = $tipo
$ws.Cells.Item(94, 7).Value = $productoId
$ws.Cells.Item(94, 8).Value = $producto
$ws.Cells.Item(94, 9).Value = $categoriaId
$ws.Cells.Item(94, 10).Value = $categoria
$ws.Cells.Item(94, 11).Value = "Black Amber"
$ws.Cells.Item(94, 12).Value = "Primera"
$ws.Cells.Item(94, 13).Value = 300
$ws.Cells.Item(94, 14).Value = 15000
$ws.Cells.Item(94, 15).Value = 16000
$ws.Cells.Item(94, 16).Value = 15500
$ws.Cells.Item(94, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(94, 18).Value = $origenOHiggins
$ws.Cells.Item(94, 19).Value = 1033
$ws.Cells.Item(94, 20).Value = 15

# New row 95: Black Amber, Segunda
$ws.Cells.Item(95, 1).Value = $mercadoId
$ws.Cells.Item(95, 2).Value = $mercado
$ws.Cells.Item(95, 3).Value = $region
$ws.Cells.Item(95, 4).Value = 44603
$ws.Cells.Item(95, 5).Value = $codreg
$ws.Cells.Item(95, 6).Value = $tipo
$ws.Cells.Item(95, 7).Value = $productoId
$ws.Cells.Item(95, 8).Value = $producto
$ws.Cells.Item(95, 9).Value = $categoriaId
$ws.Cells.Item(95, 10).Value = $categoria
$ws.Cells.Item(95, 11).Value = "Black Amber"
$ws.Cells.Item(95, 12).Value = "Segunda"
$ws.Cells.Item(95, 13).Value = 150
$ws.Cells.Item(95, 14).Value = 13000
$ws.Cells.Item(95, 15).Value = 13000
$ws.Cells.Item(95, 16).Value = 13000
$ws.Cells.Item(95, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(95, 18).Value = $origenOHiggins
$ws.Cells.Item(95, 19).Value = 867
$ws.Cells.Item(95, 20).Value = 15

# New row 96: Lemon, Primera
$ws.Cells.Item(96, 1).Value = $mercadoId
$ws.Cells.Item(96, 2).Value = $mercado
$ws.Cells.Item(96, 3).Value = $region
$ws.Cells.Item(96, 4).Value = 44603
$ws.Cells.Item(96, 5).Value = $codreg
$ws.Cells.Item(96, 6).Value = $tipo
$ws.Cells.Item(96, 7).Value = $productoId
$ws.Cells.Item(96, 8).Value = $producto
$ws.Cells.Item(96, 9).Value = $categoriaId
$ws.Cells.Item(96, 10).Value = $categoria
$ws.Cells.Item(96, 11).Value = "Lemon"
$ws.Cells.Item(96, 12).Value = "Primera"
$ws.Cells.Item(96, 13).Value = 300
$ws.Cells.Item(96, 14).Value = 15000
$ws.Cells.Item(96, 15).Value = 16000
$ws.Cells.Item(96, 16).Value = 15500
$ws.Cells.Item(96, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(96, 18).Value = $origenOHiggins
$ws.Cells.Item(96, 19).Value = 1033
$ws.Cells.Item(96, 20).Value = 15

# New row 97: Lemon, Segunda
$ws.Cells.Item(97, 1).Value = $mercadoId
$ws.Cells.Item(97, 2).Value = $mercado
$ws.Cells.Item(97, 3).Value = $region
$ws.Cells.Item(97, 4).Value = 44603
$ws.Cells.Item(97, 5).Value = $codreg
$ws.Cells.Item(97, 6).Value = $tipo
$ws.Cells.Item(97, 7).Value = $productoId
$ws.Cells.Item(97, 8).Value = $producto
$ws.Cells.Item(97, 9).Value = $categoriaId
$ws.Cells.Item(97, 10).Value = $categoria
$ws.Cells.Item(97, 11).Value = "Lemon"
$ws.Cells.Item(97, 12).Value = "Segunda"
$ws.Cells.Item(97, 13).Value = 150
$ws.Cells.Item(97, 14).Value = 13000
$ws.Cells.Item(97, 15).Value = 13000
$ws.Cells.Item(97, 16).Value = 13000
$ws.Cells.Item(97, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(97, 18).Value = $origenOHiggins
$ws.Cells.Item(97, 19).Value = 867
$ws.Cells.Item(97, 20).Value = 15

# Ensure the new rows share the same date cell style (numFmtId for dates) as the row above.
$ws.Range("D93").Copy()
$ws.Range("D94:D97").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
